$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows right above the current row 881 (shifts existing rows 881-930 down to 884-933)
$ws.Rows.Item(881).Insert()
$ws.Rows.Item(881).Insert()
$ws.Rows.Item(881).Insert()

# New row 881 - Agricola del Norte S.A. de Arica, Tomate, Larga vida, Primera - new weekly entry
$ws.Range("A881").Value = 1
$ws.Range("B881").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C881").Value = "Arica y Parinacota"
$ws.Range("D881").Value = 44746
$ws.Range("E881").Value = 15
$ws.Range("F881").Value = 100112020
$ws.Range("G881").Value = "Tomate"
$ws.Range("H881").Value = "Larga vida"
$ws.Range("I881").Value = "Primera"
$ws.Range("J881").Value = 300
$ws.Range("K881").Value = 3000
$ws.Range("L881").Value = 3500
$ws.Range("M881").Value = 3250
$ws.Range("N881").Value = "$/caja 10 kilos"
$ws.Range("O881").Value = "Región de Arica y Parinacota"
$ws.Range("P881").Value = 325
$ws.Range("Q881").Value = 10
$ws.Range("R881").Value = "Hortaliza"

# New row 882 - Segunda
$ws.Range("A882").Value = 1
$ws.Range("B882").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C882").Value = "Arica y Parinacota"
$ws.Range("D882").Value = 44746
$ws.Range("E882").Value = 15
$ws.Range("F882").Value = 100112020
$ws.Range("G882").Value = "Tomate"
$ws.Range("H882").Value = "Larga vida"
$ws.Range("I882").Value = "Segunda"
$ws.Range("J882").Value = 350
$ws.Range("K882").Value = 2500
$ws.Range("L882").Value = 3000
$ws.Range("M882").Value = 2750
$ws.Range("N882").Value = "$/caja 10 kilos"
$ws.Range("O882").Value = "Región de Arica y Parinacota"
$ws.Range("P882").Value = 275
$ws.Range("Q882").Value = 10
$ws.Range("R882").Value = "Hortaliza"

# New row 883 - Tercera
$ws.Range("A883").Value = 1
$ws.Range("B883").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C883").Value = "Arica y Parinacota"
$ws.Range("D883").Value = 44746
$ws.Range("E883").Value = 15
$ws.Range("F883").Value = 100112020
$ws.Range("G883").Value = "Tomate"
$ws.Range("H883").Value = "Larga vida"
$ws.Range("I883").Value = "Tercera"
$ws.Range("J883").Value = 350
$ws.Range("K883").Value = 2000
$ws.Range("L883").Value = 2500
$ws.Range("M883").Value = 2250
$ws.Range("N883").Value = "$/caja 10 kilos"
$ws.Range("O883").Value = "Región de Arica y Parinacota"
$ws.Range("P883").Value = 225
$ws.Range("Q883").Value = 10
$ws.Range("R883").Value = "Hortaliza"
